$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values taken from the authoritative diff: columns B (r_value SUM),
# C (Word Count), F (Transcribed), G (Weight, mirrors B) and H (Result) are
# updated for data rows 2-41; D (NonTranscribed) and E (NBiphones) are
# untouched by the commit.
$rows = @(
    @{Row=2; B=127.0108868979071; C=6136; F=6021; H=0.04284750337381917},
    @{Row=3; B=-61.24612018567306; C=732; F=715; H=-0.1452380952380952},
    @{Row=4; B=-212.3418476210031; C=3904; F=3766; H=-0.1130063965884861},
    @{Row=5; B=-83.25904031629506; C=1444; F=1409; H=-0.1034912718204489},
    @{Row=6; B=-381.8469510135291; C=3744; F=3660; H=-0.1753336401288541},
    @{Row=7; B=-226.9751793531174; C=3392; F=3358; H=-0.1197668256491786},
    @{Row=8; B=-156.8876156948677; C=3296; F=3207; H=-0.08436992969172526},
    @{Row=9; B=-82.32936737153875; C=3924; F=3827; H=-0.04039408866995074},
    @{Row=10; B=37.18195107458128; C=1632; F=1611; H=0.03775510204081633},
    @{Row=11; B=-149.662743086256; C=2608; F=2574; H=-0.1005398110661269},
    @{Row=12; B=-260.2462975575222; C=1844; F=1824; H=-0.2561576354679803},
    @{Row=13; B=-125.4879013688354; C=888; F=873; H=-0.2642706131078224},
    @{Row=14; B=45.489393271398; C=1464; F=1444; H=0.05528255528255528},
    @{Row=15; B=-720.0590425069661; C=7348; F=7176; H=-0.1748421563865954},
    @{Row=16; B=47.70895354011265; C=2560; F=2519; H=0.03440702781844802},
    @{Row=17; B=-594.9757865125223; C=3896; F=3845; H=-0.2884895580378825},
    @{Row=18; B=-32.93947489898112; C=976; F=964; H=-0.06557377049180328},
    @{Row=19; B=-450.8703413625891; C=4636; F=4550; H=-0.1861042183622829},
    @{Row=20; B=-112.7130100074113; C=1188; F=1161; H=-0.2021660649819494},
    @{Row=21; B=-50.26298722493762; C=2252; F=2229; H=-0.0370919881305638},
    @{Row=22; B=-109.5953780685138; C=2624; F=2592; H=-0.07185234014502308},
    @{Row=23; B=19.07818271353684; C=1228; F=1202; H=0.02303030303030303},
    @{Row=24; B=-141.1988355505335; C=1956; F=1911; H=-0.1180904522613065},
    @{Row=25; B=20.9436434891197; C=1408; F=1396; H=0.02680965147453083},
    @{Row=26; B=-282.5282828522926; C=2940; F=2872; H=-0.1934156378600823},
    @{Row=27; B=-59.25921405929971; C=1244; F=1208; H=-0.07919463087248323},
    @{Row=28; B=-226.8913927195349; C=1628; F=1580; H=-0.2807453416149068},
    @{Row=29; B=-448.4560768004027; C=2496; F=2415; H=-0.3592622293504411},
    @{Row=30; B=191.5671265455729; C=6580; F=6469; H=0.04612412460758271},
    @{Row=31; B=-85.95745261293301; C=784; F=760; H=-0.2048192771084337},
    @{Row=32; B=-102.6527857225576; C=3068; F=2997; H=-0.04973183812774257},
    @{Row=33; B=-128.8681815330119; C=1716; F=1651; H=-0.1559074299634592},
    @{Row=34; B=-363.6641642799185; C=2900; F=2865; H=-0.2178871548619448},
    @{Row=35; B=-124.3523736193784; C=1980; F=1921; H=-0.1150278293135436},
    @{Row=36; B=58.44100642036419; C=1692; F=1641; H=0.04869857262804366},
    @{Row=37; B=48.02254751488483; C=2892; F=2829; H=0.03104786545924968},
    @{Row=38; B=-88.97232123756133; C=736; F=726; H=-0.2156862745098039},
    @{Row=39; B=-504.1183246364689; C=3320; F=3266; H=-0.2807799442896936},
    @{Row=40; B=-362.1260193438251; C=3424; F=3367; H=-0.1776251226692836},
    @{Row=41; B=-65.5626212472756; C=2140; F=2083; H=-0.05617977528089887}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B   # B: r_value SUM
    $ws.Cells.Item($r.Row, 3).Value = $r.C   # C: Word Count
    $ws.Cells.Item($r.Row, 6).Value = $r.F   # F: Transcribed
    $ws.Cells.Item($r.Row, 7).Value = $r.B   # G: Weight (mirrors B)
    $ws.Cells.Item($r.Row, 8).Value = $r.H   # H: Result
}
